$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: new simulation entry (ID "011", recovery rate lowered to 0.7) ---
$ws.Range("B12").Value = "1 / 1.8"
$ws.Range("C12").Value = 0.7
$ws.Range("D12").Value = 1.8
$ws.Range("E12").Value = 1000
$ws.Range("F12").Value = "children, adolescents, adults, elderly"
$ws.Range("G12").Value = "open, close"
$ws.Range("H12").Value = 43
$ws.Range("I12").Value = 0.5
$ws.Range("J12").Value = $false
$ws.Range("K12").Value = 1000000

# --- Row 13: new simulation entry (ID "012", recovery rate 0.8) ---
$ws.Range("B13").Value = "1 / 1.8"
$ws.Range("C13").Value = 0.8
$ws.Range("D13").Value = 1.8
$ws.Range("E13").Value = 1000
$ws.Range("F13").Value = "children, adolescents, adults, elderly"
$ws.Range("G13").Value = "open, close"
$ws.Range("H13").Value = 43
$ws.Range("I13").Value = 0.5
$ws.Range("J13").Value = $false
$ws.Range("K13").Value = 1000000

# Match the highlighted "budget" cell formatting used by the other rows
$ws.Range("E11").Copy()
$ws.Range("E12").PasteSpecial(-4122)
$ws.Range("E11").Copy()
$ws.Range("E13").PasteSpecial(-4122)

# Update the saved selection to reflect where the user left off
$ws.Range("J28").Select() | Out-Null
